$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook stores Price/Volume columns as plain text (not numbers),
# even though many Price values look numeric. Excel COM auto-converts
# numeric-looking strings assigned via .Value into real numbers, which
# changes both the stored type and introduces float rounding artifacts.
# To keep them as literal text (matching the original formatting,
# e.g. "313.59" and not "313.58999999999997"), we temporarily force the
# cell to Text format, assign the literal string, then strip the
# number-format override back off so the cell keeps the workbook's
# original (default) style.

# Row 2
$ws.Range("D2").Value = "44.442.08"
$ws.Range("E2").Value = "  +3.47%  "

# Row 3
$ws.Range("D3").Value = "2.420.51"
$ws.Range("E3").Value = "  +2.41%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.70%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.04%  "

# Row 7
$ws.Range("E7").Value = "  +1.86%  "

# Row 8
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.512"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.45%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.00%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0798"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.76%  "

# Row 12
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.70%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.125"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.72%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.05%  "

# Row 15
$ws.Range("D15").Value = "2.799.71"
$ws.Range("E15").Value = "  +2.55%  "

# Row 16
$ws.Range("D16").Value = "2.433.22"
$ws.Range("E16").Value = "  +2.23%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.832"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.10%  "

# Row 18
$ws.Range("D18").Value = "44.320.95"
$ws.Range("E18").Value = "  +3.28%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.68%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.67%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0918"
$ws.Range("E21").Value = "  +3.95%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.04%  "

# Row 24
$ws.Range("E24").Value = "  +5.75%  "

# Row 25
$ws.Range("E25").Value = "  +1.40%  "

# Row 26
$ws.Range("E26").Value = "  -0.06%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.62%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.13%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.58%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.91%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.97%  "

# Row 32
$ws.Range("E32").Value = "  +18.59%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.32"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.37%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.78%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0770"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.99%  "

# Row 36
$ws.Range("E36").Value = "  +0.24%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.13%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.95%  "

# Row 39
$ws.Range("E39").Value = "  +2.65%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "122.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.78%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.34%  "

# Row 42
$ws.Range("E42").Value = "  +1.28%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.15%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0288"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.66%  "

# Row 45
$ws.Range("D45").Value = "1.949.27"
$ws.Range("E45").Value = "  +0.98%  "

# Row 46
$ws.Range("E46").Value = "  +2.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.82%  "

# Row 48
$ws.Range("E48").Value = "  +3.43%  "

# Row 49
$ws.Range("E49").Value = "  +9.95%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.59%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.18%  "
